# Add ErrorCode service implementation
# Fills in the "id" column (A) for rows 4 through 23 with sequential
# error code values (1002 .. 1021), and updates the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate column A (id) for rows 4..23 with values 1002..1021
$startRow = 4
$startValue = 1002
for ($row = $startRow; $row -le 23; $row++) {
    $value = $startValue + ($row - $startRow)
    $ws.Cells.Item($row, 1).Value = $value
}

# Update the active selection to match the authored change
$ws.Range("E19").Select()
